$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.201.48'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.563.60'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +4.92%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.65'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.85'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.563.37'
$ws.Range("D7").ClearFormats()
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.502'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.13%  '
$ws.Range("E10").Value = '  +2.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.97'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("E12").Value = '  +2.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.167.73'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.94%  '
$ws.Range("E14").Value = '  +4.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.11'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.563.00'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.311.14'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.37'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +10.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.20'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.92'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.88'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.48%  '
$ws.Range("E23").Value = '  +6.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.95'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.704.84'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.79%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000119'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +9.88%  '
$ws.Range("E28").Value = '  +4.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.01'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.12'
$ws.Range("D30").ClearFormats()
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("E32").Value = '  +1.44%  '
$ws.Range("E33").Value = '  +5.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.557.99'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.153'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.98%  '
$ws.Range("E37").Value = '  +4.72%  '
$ws.Range("E38").Value = '  +5.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.62'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.27%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '172.84'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0855'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.22'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.896'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.15'
$ws.Range("D46").ClearFormats()
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.85'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.35'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.87%  '
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.43'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +16.16%  '
